$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Hydrogen / Iron & steel demand value
$ws.Range("B3").Value = 3139.000735912728

# Update Other / Non-metallic minerals value
$ws.Range("D8").Value = 3150.187038115713
